# C5-PowerPoint.pptx edit script
# 1) Re-point the sources-of-finance table (slide 6) at the built-in
#    PowerPoint table style {8FC63198-2A05-4FB6-9A85-97D910E604CD}
#    instead of the custom "Table_0" style.
# 2) Re-theme the deck: swap the "Integral" theme's colour scheme for the
#    stock "Office Theme" colour scheme (the font scheme and format scheme
#    are already identical between the two themes, so only the 12 scheme
#    colours need to change).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$s6 = $p.Slides.Item(6)
for ($i = 1; $i -le $s6.Shapes.Count; $i++) {
    $shp = $s6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{8FC63198-2A05-4FB6-9A85-97D910E604CD}")
    }
}

# --- 2. Theme colours -------------------------------------------------
# msoThemeColorDark1..msoThemeColorFollowedHyperlink, 1-based, RRGGBB
$officeColors = @(
    "000000",  # Dark 1
    "FFFFFF",  # Light 1
    "44546A",  # Dark 2
    "E7E6E6",  # Light 2
    "5B9BD5",  # Accent 1
    "ED7D31",  # Accent 2
    "A5A5A5",  # Accent 3
    "FFC000",  # Accent 4
    "4472C4",  # Accent 5
    "70AD47",  # Accent 6
    "0563C1",  # Hyperlink
    "954F72"   # Followed Hyperlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # PowerPoint's ColorFormat.RGB is a VBA RGB() long: 0x00BBGGRR
    $tcs.Item($i).RGB = ($b * 65536) + ($g * 256) + $r
}
